$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.064198653302102
$ws.Range("D2").Value = 1.06175269416624
$ws.Range("E2").Value = 1.069076882600384
$ws.Range("F2").Value = 1.078871629810596
$ws.Range("I2").Value = 1.055988313694664
$ws.Range("J2").Value = 1.069160123005274
$ws.Range("K2").Value = 1.064475900687661
$ws.Range("L2").Value = 1.071780355622972
$ws.Range("M2").Value = 1.081549165381923
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.065430994089894
$ws.Range("D3").Value = 1.062697770233633
$ws.Range("E3").Value = 1.070189976582276
$ws.Range("F3").Value = 1.080090757819483
$ws.Range("I3").Value = 1.05642783024316
$ws.Range("J3").Value = 1.070046242116565
$ws.Range("K3").Value = 1.065235330747638
$ws.Range("L3").Value = 1.072708813494198
$ws.Range("M3").Value = 1.082585280450773
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.066227995820213
$ws.Range("D4").Value = 1.063308840469844
$ws.Range("E4").Value = 1.07091003792658
$ws.Range("F4").Value = 1.080879649755082
$ws.Range("I4").Value = 1.056710730055846
$ws.Range("J4").Value = 1.070618669374011
$ws.Range("K4").Value = 1.065725650571288
$ws.Range("L4").Value = 1.073308810367187
$ws.Range("M4").Value = 1.083255176655955
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.066562960649962
$ws.Range("D5").Value = 1.063565626036134
$ws.Range("E5").Value = 1.071212708961711
$ws.Range("F5").Value = 1.081211309921249
$ws.Range("I5").Value = 1.056829303615877
$ws.Range("J5").Value = 1.070859091538858
$ws.Range("K5").Value = 1.065931523095783
$ws.Range("L5").Value = 1.073560864351585
$ws.Range("M5").Value = 1.083536673403854
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.06661919723251
$ws.Range("D6").Value = 1.063608735139855
$ws.Range("E6").Value = 1.071263526297081
$ws.Range("F6").Value = 1.081266997758558
$ws.Range("I6").Value = 1.056849191696732
$ws.Range("J6").Value = 1.070899446252187
$ws.Range("K6").Value = 1.065966074906923
$ws.Range("L6").Value = 1.073603174535057
$ws.Range("M6").Value = 1.083583930507926
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.066232472005883
$ws.Range("D7").Value = 1.063312272075578
$ws.Range("E7").Value = 1.070914082397872
$ws.Range("F7").Value = 1.080884081373031
$ws.Range("I7").Value = 1.056712315845361
$ws.Range("J7").Value = 1.070621882793983
$ws.Range("K7").Value = 1.065728402462086
$ws.Range("L7").Value = 1.073312179048835
$ws.Range("M7").Value = 1.083258938528793
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.064615213917316
$ws.Range("D8").Value = 1.062072182364862
$ws.Range("E8").Value = 1.069453097115555
$ws.Range("F8").Value = 1.079283633699152
$ws.Range("I8").Value = 1.056137160777981
$ws.Range("J8").Value = 1.069459788893502
$ws.Range("K8").Value = 1.064732778404907
$ws.Range("L8").Value = 1.072094294141567
$ws.Range("M8").Value = 1.081899437965073
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.06176219230467
$ws.Range("D9").Value = 1.059883435266897
$ws.Range("E9").Value = 1.066877170374092
$ws.Range("F9").Value = 1.076463628233203
$ws.Range("I9").Value = 1.055112163930702
$ws.Range("J9").Value = 1.067404682556591
$ws.Range("K9").Value = 1.062970021233544
$ws.Range("L9").Value = 1.069942206946032
$ws.Range("M9").Value = 1.079499622537262
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.059857875839257
$ws.Range("D10").Value = 1.058421812081614
$ws.Range("E10").Value = 1.065158790685911
$ws.Range("F10").Value = 1.074583640562618
$ws.Range("I10").Value = 1.054421048199816
$ws.Range("J10").Value = 1.066029578861577
$ws.Range("K10").Value = 1.061789167547
$ws.Range("L10").Value = 1.068503344111924
$ws.Range("M10").Value = 1.077896820993715
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.059032706522832
$ws.Range("D11").Value = 1.057788313359613
$ws.Range("E11").Value = 1.064414431103092
$ws.Range("F11").Value = 1.073769563612398
$ws.Range("I11").Value = 1.054119928966758
$ws.Range("J11").Value = 1.065432928899488
$ws.Range("K11").Value = 1.061276479869947
$ws.Range("L11").Value = 1.067879298075046
$ws.Range("M11").Value = 1.077202073679487
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.058726110640814
$ws.Range("D12").Value = 1.057552911116604
$ws.Range("E12").Value = 1.064137897236372
$ws.Range("F12").Value = 1.073467172417559
$ws.Range("I12").Value = 1.054007799083725
$ws.Range("J12").Value = 1.065211121181613
$ws.Range("K12").Value = 1.061085837401707
$ws.Range("L12").Value = 1.067647346103435
$ws.Range("M12").Value = 1.076943902803432
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.05879188067544
$ws.Range("D13").Value = 1.057603409912321
$ws.Range("E13").Value = 1.064197216784531
$ws.Range("F13").Value = 1.073532036698312
$ws.Range("I13").Value = 1.054031864048638
$ws.Range("J13").Value = 1.065258708108627
$ws.Range("K13").Value = 1.061126740270619
$ws.Range("L13").Value = 1.067697107565046
$ws.Range("M13").Value = 1.076999286379454
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.059007365090339
$ws.Range("D14").Value = 1.057768856834888
$ws.Range("E14").Value = 1.064391573668884
$ws.Range("F14").Value = 1.073744568013187
$ws.Range("I14").Value = 1.054110666010182
$ws.Range("J14").Value = 1.065414597996647
$ws.Range("K14").Value = 1.061260725544839
$ws.Range("L14").Value = 1.067860127994866
$ws.Range("M14").Value = 1.077180735466391
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.059140119969148
$ws.Range("D15").Value = 1.057870781826544
$ws.Range("E15").Value = 1.064511317259957
$ws.Range("F15").Value = 1.07387551460549
$ws.Range("I15").Value = 1.054159181258382
$ws.Range("J15").Value = 1.06551062228277
$ws.Range("K15").Value = 1.061343250759682
$ws.Range("L15").Value = 1.06796054986628
$ws.Range("M15").Value = 1.077292517494769
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.059912626974161
$ws.Range("D16").Value = 1.058463842402174
$ws.Range("E16").Value = 1.065208185100226
$ws.Range("F16").Value = 1.074637667344712
$ws.Range("I16").Value = 1.054440993175614
$ws.Range("J16").Value = 1.066069150666267
$ws.Range("K16").Value = 1.061823163935976
$ws.Range("L16").Value = 1.0685447385724
$ws.Range("M16").Value = 1.0779429136617
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.060397040772177
$ws.Range("D17").Value = 1.058835690246593
$ws.Range("E17").Value = 1.065645232855507
$ws.Range("F17").Value = 1.075115736014019
$ws.Range("I17").Value = 1.054617267167105
$ws.Range("J17").Value = 1.066419172520121
$ws.Range("K17").Value = 1.062123832829888
$ws.Range("L17").Value = 1.068910913248004
$ws.Range("M17").Value = 1.078350694789242
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.060679534639768
$ws.Range("D18").Value = 1.059052524414766
$ws.Range("E18").Value = 1.065900127590293
$ws.Range("F18").Value = 1.075394582652771
$ws.Range("I18").Value = 1.054719905266353
$ws.Range("J18").Value = 1.066623216668644
$ws.Range("K18").Value = 1.062299075808165
$ws.Range("L18").Value = 1.06912439952358
$ws.Range("M18").Value = 1.078588476939238
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.060775848253015
$ws.Range("D19").Value = 1.059126449404728
$ws.Range("E19").Value = 1.065987035448774
$ws.Range("F19").Value = 1.075489661756233
$ws.Range("I19").Value = 1.054754871786513
$ws.Range("J19").Value = 1.066692770547775
$ws.Range("K19").Value = 1.062358806775725
$ws.Range("L19").Value = 1.069197176394003
$ws.Range("M19").Value = 1.078669542750175
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.060345073623348
$ws.Range("D20").Value = 1.05879580054339
$ws.Range("E20").Value = 1.065598344661536
$ws.Range("F20").Value = 1.075064444069987
$ws.Range("I20").Value = 1.054598373215116
$ws.Range("J20").Value = 1.066381630674465
$ws.Range("K20").Value = 1.062091587579084
$ws.Range("L20").Value = 1.068871636241179
$ws.Range("M20").Value = 1.078306950957051
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.058943912843901
$ws.Range("D21").Value = 1.057720139408692
$ws.Range("E21").Value = 1.064334341685086
$ws.Range("F21").Value = 1.073681983043345
$ws.Range("I21").Value = 1.0540874685568
$ws.Range("J21").Value = 1.06536869744223
$ws.Range("K21").Value = 1.061221275960037
$ws.Range("L21").Value = 1.067812126793813
$ws.Range("M21").Value = 1.077127306307756
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.058062417993863
$ws.Range("D22").Value = 1.057043291978456
$ws.Range("E22").Value = 1.063539348153076
$ws.Range("F22").Value = 1.072812734631215
$ws.Range("I22").Value = 1.053764617388112
$ws.Range("J22").Value = 1.064730752749945
$ws.Range("K22").Value = 1.06067287554977
$ws.Range("L22").Value = 1.067145082171477
$ws.Range("M22").Value = 1.076384974873138
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.058529766268749
$ws.Range("D23").Value = 1.057402152982424
$ws.Range("E23").Value = 1.063960815035654
$ws.Range("F23").Value = 1.073273544243703
$ws.Range("I23").Value = 1.053935921317681
$ws.Range("J23").Value = 1.065069041685889
$ws.Range("K23").Value = 1.060963707347039
$ws.Range("L23").Value = 1.067498780050008
$ws.Range("M23").Value = 1.076778560360275
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.060368555515566
$ws.Range("D24").Value = 1.058813825163714
$ws.Range("E24").Value = 1.06561953150414
$ws.Range("F24").Value = 1.075087620700129
$ws.Range("I24").Value = 1.05460691113349
$ws.Range("J24").Value = 1.066398594582165
$ws.Range("K24").Value = 1.062106158228138
$ws.Range("L24").Value = 1.068889384130315
$ws.Range("M24").Value = 1.078326717127224
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.062500161785686
$ws.Range("D25").Value = 1.060449706932376
$ws.Range("E25").Value = 1.067543295058515
$ws.Range("F25").Value = 1.077192656323183
$ws.Range("I25").Value = 1.055378518403893
$ws.Range("J25").Value = 1.067936857101575
$ws.Range("K25").Value = 1.063426732263352
$ws.Range("L25").Value = 1.070499295699918
$ws.Range("M25").Value = 1.08012054051313
